$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("3G SNC DI PAOLO GRANELLI &amp; CO.", "amministrazione@salsanatura.com"),
    @("3G SNC DI PAOLO GRANELLI &amp; CO.", "amministrazione@salsanatura.com"),
    @("A RICCHIGIA SRL", "aricchigia@gmail.com"),
    @("3G SNC DI PAOLO GRANELLI &amp; CO.", "amministrazione@salsanatura.com"),
    @("A RICCHIGIA SRL", "aricchigia@gmail.com"),
    @("A. DARBO AG", "katrin.widauer@darbo.at"),
    @("A. GANDOLA &amp; C. SPA", "gandola@gandola.it"),
    @("A. GANDOLA &amp; C. SPA", "gandola@gandola.it"),
    @("A. LOACKER SPA", "marketing@loacker.com")
)

$row = 4
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
